$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
$ws.Range("F1").Value = "Investment"
$ws.Range("K1").Value = "Land Saving"
$ws.Range("L1").Value = "Workforce Saving"
$ws.Range("M1").Value = "Water Investment"
$ws.Range("N1").Value = "Emission Investment"
$ws.Range("O1").Value = "Land Investment"
$ws.Range("P1").Value = "Workforce Investment"

# New header cells O1/P1 need the same bold/centered/bordered style as the
# rest of row 1 - copy formatting from the neighboring N1 cell.
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

# --- Row 2 units ---
$ws.Range("C2").Value = "M kSh/FU"
$ws.Range("E2").Value = "M kSh"
$ws.Range("F2").Value = "M kSh"
$ws.Range("G2").Value = "M kSh/FU"
$ws.Range("H2").Value = "years"
$ws.Range("J2").Value = "kton/FU"
$ws.Range("K2").Value = "M kSh/FU"
$ws.Range("L2").Value = "M kSh/FU"
$ws.Range("M2").Value = "m3/FU"
$ws.Range("N2").Value = "kton/FU"
$ws.Range("O2").Value = "M kSh/FU"
$ws.Range("P2").Value = "M kSh/FU"

# New unit cells O2/P2 likewise need matching formatting from N2.
$ws.Range("N2").Copy()
$ws.Range("O2:P2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 4 identifiers ---
$ws.Range("B4").Value = "17 plants"

# --- Row 4 data values ---
$ws.Range("F4").Value = 5421.299537870102
$ws.Range("G4").Value = 3253.734504770488
$ws.Range("H4").Value = 1.666177596826546
$ws.Range("I4").Value = 9.686076806276105
$ws.Range("J4").Value = 58.23311956987891
$ws.Range("K4").Value = 3.371160993818194
$ws.Range("L4").Value = 75.68173810420558
$ws.Range("M4").Value = 24.33238917902054
$ws.Range("N4").Value = 14.82826949442824
$ws.Range("O4").Value = 6.616823517833836
$ws.Range("P4").Value = 278.0264926441014
